$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.084220245949388
$ws.Range("D2").Value = 1.074005857132056
$ws.Range("E2").Value = 1.098084468359911
$ws.Range("F2").Value = 1.105424445632466
$ws.Range("I2").Value = 1.045971304505695
$ws.Range("J2").Value = 1.089081098627202
$ws.Range("K2").Value = 1.07669621384224
$ws.Range("L2").Value = 1.100712575745193
$ws.Range("M2").Value = 1.108034158292004
$ws.Range("N2").Value = 1.090627718706601
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.086755652850841
$ws.Range("D3").Value = 1.075957790802158
$ws.Range("E3").Value = 1.100662152279938
$ws.Range("F3").Value = 1.108153581089982
$ws.Range("I3").Value = 1.046613614207576
$ws.Range("J3").Value = 1.091272028890468
$ws.Range("K3").Value = 1.078462404626961
$ws.Range("L3").Value = 1.103107657053695
$ws.Range("M3").Value = 1.110581719848995
$ws.Range("N3").Value = 1.092821760342144
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.088386999326441
$ws.Range("D4").Value = 1.077212485887669
$ws.Range("E4").Value = 1.102321109918172
$ws.Range("F4").Value = 1.109910540768271
$ws.Range("I4").Value = 1.047024026027149
$ws.Range("J4").Value = 1.092680453815718
$ws.Range("K4").Value = 1.079596509649723
$ws.Range("L4").Value = 1.104648167535615
$ws.Range("M4").Value = 1.112220921855374
$ws.Range("N4").Value = 1.094232185392336
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.08907066281824
$ws.Range("D5").Value = 1.077738006243968
$ws.Range("E5").Value = 1.103016442064054
$ws.Range("F5").Value = 1.1106470764706
$ws.Range("I5").Value = 1.04719533338457
$ws.Range("J5").Value = 1.093270389416624
$ws.Range("K5").Value = 1.0800712358212
$ws.Range("L5").Value = 1.105293632205171
$ws.Range("M5").Value = 1.112907885502369
$ws.Range("N5").Value = 1.094822958769462
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.089185328178443
$ws.Range("D6").Value = 1.077826130091266
$ws.Range("E6").Value = 1.103133070210766
$ws.Range("F6").Value = 1.110770623162653
$ws.Range("I6").Value = 1.047224025002471
$ws.Range("J6").Value = 1.093369316583449
$ws.Range("K6").Value = 1.080150825293064
$ws.Range("L6").Value = 1.105401883068619
$ws.Range("M6").Value = 1.113023104835601
$ws.Range("N6").Value = 1.094922026424211
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.08839614286368
$ws.Range("D7").Value = 1.077219515535755
$ws.Range("E7").Value = 1.102330409133349
$ws.Range("F7").Value = 1.109920390529876
$ws.Range("I7").Value = 1.047026319854665
$ws.Range("J7").Value = 1.092688345010607
$ws.Range("K7").Value = 1.079602860971498
$ws.Range("L7").Value = 1.104656800706663
$ws.Range("M7").Value = 1.112230109498614
$ws.Range("N7").Value = 1.094240087793628
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.085079043651648
$ws.Range("D8").Value = 1.074667274090745
$ws.Range("E8").Value = 1.09895750086498
$ws.Range("F8").Value = 1.106348663041825
$ws.Range("I8").Value = 1.046189463897955
$ws.Range("J8").Value = 1.089823479551548
$ws.Range("K8").Value = 1.077294940216452
$ws.Range("L8").Value = 1.101523955523216
$ws.Range("M8").Value = 1.108897065607461
$ws.Range("N8").Value = 1.091371153897012
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.079160665506952
$ws.Range("D9").Value = 1.07010412259846
$ws.Range("E9").Value = 1.09294274081736
$ws.Range("F9").Value = 1.099983450358093
$ws.Range("I9").Value = 1.044674199271184
$ws.Range("J9").Value = 1.084702115333644
$ws.Range("K9").Value = 1.073159319734687
$ws.Range("L9").Value = 1.095930106822572
$ws.Range("M9").Value = 1.102950522652918
$ws.Range("N9").Value = 1.086242516754489
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.075162268683193
$ws.Range("D10").Value = 1.067015062993851
$ws.Range("E10").Value = 1.088881417718334
$ws.Range("F10").Value = 1.095688237240138
$ws.Range("I10").Value = 1.04363564971705
$ws.Range("J10").Value = 1.081235551925881
$ws.Range("K10").Value = 1.070353389230157
$ws.Range("L10").Value = 1.092148140040682
$ws.Range("M10").Value = 1.09893330806061
$ws.Range("N10").Value = 1.082771030429067
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.073417581753719
$ws.Range("D11").Value = 1.065665708204038
$ws.Range("E11").Value = 1.08710980037483
$ws.Range("F11").Value = 1.093815241227543
$ws.Range("I11").Value = 1.043178974940675
$ws.Range("J11").Value = 1.079721363786072
$ws.Range("K11").Value = 1.06912621110758
$ws.Range("L11").Value = 1.090497228145679
$ws.Range("M11").Value = 1.097180464093219
$ws.Range("N11").Value = 1.081254691968419
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.072767450212553
$ws.Range("D12").Value = 1.065162675076506
$ws.Range("E12").Value = 1.086449713520124
$ws.Range("F12").Value = 1.093117478799204
$ws.Range("I12").Value = 1.043008276386073
$ws.Range("J12").Value = 1.079156889572312
$ws.Range("K12").Value = 1.068668499086374
$ws.Range("L12").Value = 1.089881941240354
$ws.Range("M12").Value = 1.096527301428993
$ws.Range("N12").Value = 1.08068941613653
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.07290700078457
$ws.Range("D13").Value = 1.065270660753752
$ws.Range("E13").Value = 1.086591397372421
$ws.Range("F13").Value = 1.093267245109246
$ws.Range("I13").Value = 1.043044940534834
$ws.Range("J13").Value = 1.079278064423757
$ws.Range("K13").Value = 1.068766765942819
$ws.Range("L13").Value = 1.090014016901058
$ws.Range("M13").Value = 1.096667502241668
$ws.Range("N13").Value = 1.080810763070163
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.073363884465927
$ws.Range("D14").Value = 1.065624164828705
$ws.Range("E14").Value = 1.08705527920056
$ws.Range("F14").Value = 1.093757606187693
$ws.Range("I14").Value = 1.043164886880649
$ws.Range("J14").Value = 1.079674746100234
$ws.Range("K14").Value = 1.069088415214388
$ws.Range("L14").Value = 1.090446410852662
$ws.Range("M14").Value = 1.097126516276538
$ws.Range("N14").Value = 1.081208008080121
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.073645108196831
$ws.Range("D15").Value = 1.065841727193188
$ws.Range("E15").Value = 1.087340821066324
$ws.Range("F15").Value = 1.094059460164246
$ws.Range("I15").Value = 1.043238647461313
$ws.Range("J15").Value = 1.07991888290872
$ws.Range("K15").Value = 1.069286342805174
$ws.Range("L15").Value = 1.090712547403325
$ws.Range("M15").Value = 1.097409052427484
$ws.Range("N15").Value = 1.081452491590877
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.075277770439423
$ws.Range("D16").Value = 1.067104362621537
$ws.Range("E16").Value = 1.088998713455171
$ws.Range("F16").Value = 1.095812258676412
$ws.Range("I16").Value = 1.043665808948449
$ws.Range("J16").Value = 1.081335761373885
$ws.Range("K16").Value = 1.070434571744708
$ws.Range("L16").Value = 1.092257419804149
$ws.Range("M16").Value = 1.099049350994302
$ws.Range("N16").Value = 1.082871382185981
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.076298271465243
$ws.Range("D17").Value = 1.067893191055828
$ws.Range("E17").Value = 1.090035125717785
$ws.Range("F17").Value = 1.096908173028777
$ws.Range("I17").Value = 1.043931873624574
$ws.Range("J17").Value = 1.08222096916511
$ws.Range("K17").Value = 1.071151525193762
$ws.Range("L17").Value = 1.093222870922304
$ws.Range("M17").Value = 1.100074639785295
$ws.Range("N17").Value = 1.083757847073805
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.076892228768117
$ws.Range("D18").Value = 1.068352169264217
$ws.Range("E18").Value = 1.090638394355257
$ws.Range("F18").Value = 1.097546138560915
$ws.Range("I18").Value = 1.044086393021943
$ws.Range("J18").Value = 1.082736032058821
$ws.Range("K18").Value = 1.071568539243243
$ws.Range("L18").Value = 1.093784724037062
$ws.Range("M18").Value = 1.100671388853467
$ws.Range("N18").Value = 1.084273641415904
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.077094537244654
$ws.Range("D19").Value = 1.068508478575302
$ws.Range("E19").Value = 1.090843882936975
$ws.Range("F19").Value = 1.097763456354753
$ws.Range("I19").Value = 1.044138966857365
$ws.Range("J19").Value = 1.082911442705419
$ws.Range("K19").Value = 1.071710532886824
$ws.Range("L19").Value = 1.093976086641991
$ws.Range("M19").Value = 1.100874649101774
$ws.Range("N19").Value = 1.084449301165739
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.076188914646725
$ws.Range("D20").Value = 1.067808674713348
$ws.Range("E20").Value = 1.089924058609727
$ws.Range("F20").Value = 1.096790722935249
$ws.Range("I20").Value = 1.043903397023804
$ws.Range("J20").Value = 1.082126125919875
$ws.Range("K20").Value = 1.071074724549387
$ws.Range("L20").Value = 1.093119419844881
$ws.Range("M20").Value = 1.099964769372908
$ws.Range("N20").Value = 1.083662869140283
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.073229401524348
$ws.Range("D21").Value = 1.065520117445121
$ws.Range("E21").Value = 1.086918734156854
$ws.Range("F21").Value = 1.093613264219043
$ws.Range("I21").Value = 1.043129595359051
$ws.Range("J21").Value = 1.079557989983986
$ws.Range("K21").Value = 1.068993749934115
$ws.Range("L21").Value = 1.090319139155049
$ws.Range("M21").Value = 1.096991406032337
$ws.Range("N21").Value = 1.081091086156796
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.071356584293055
$ws.Range("D22").Value = 1.064070637487299
$ws.Range("E22").Value = 1.085017390424898
$ws.Range("F22").Value = 1.091603580953086
$ws.Range("I22").Value = 1.04263687859234
$ws.Range("J22").Value = 1.077931481594696
$ws.Range("K22").Value = 1.067674435170318
$ws.Range("L22").Value = 1.088546511275538
$ws.Range("M22").Value = 1.09510987351418
$ws.Range("N22").Value = 1.079462267939029
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.07235056675483
$ws.Range("D23").Value = 1.064840054852766
$ws.Range("E23").Value = 1.086026469007487
$ws.Range("F23").Value = 1.092670104196953
$ws.Range("I23").Value = 1.042898671662922
$ws.Range("J23").Value = 1.078794866117757
$ws.Range("K23").Value = 1.068374881900376
$ws.Range("L23").Value = 1.089487373304049
$ws.Range("M23").Value = 1.096108476757432
$ws.Range("N23").Value = 1.080326878567145
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.076238332248774
$ws.Range("D24").Value = 1.067846867505949
$ws.Range("E24").Value = 1.08997424892354
$ws.Range("F24").Value = 1.096843797475131
$ws.Range("I24").Value = 1.043916266448204
$ws.Range("J24").Value = 1.082168985401981
$ws.Range("K24").Value = 1.071109431075973
$ws.Range("L24").Value = 1.093166168879956
$ws.Range("M24").Value = 1.100014419048581
$ws.Range("N24").Value = 1.08370578948777
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.080699749390412
$ws.Range("D25").Value = 1.071291879284688
$ws.Range("E25").Value = 1.094506508361121
$ws.Range("F25").Value = 1.101637849632629
$ws.Range("I25").Value = 1.045070852979053
$ws.Range("J25").Value = 1.086035096076261
$ws.Range("K25").Value = 1.074236891817656
$ws.Range("L25").Value = 1.09738529173432
$ws.Range("M25").Value = 1.104496895430131
$ws.Range("N25").Value = 1.087577390482656
